$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Trim the whitespace on the indicator text in B4 (this causes Excel to
# move the edited shared string to the end of the shared-strings table).
$ws.Range("B4").Value = "12.4.1 Число сторон международных многосторонних экологических соглашений по опасным отходам и иным химических веществам, выполняющих свои обязательства и обязанности по передаче информации в соответствии с требованиями каждого соглашения"

# Update the active selection to B4 (as recorded in the saved view state).
$ws.Range("B4").Select()

# Update the window view size/position recorded in the workbook
# (maximized, starting at the screen origin).
$excel.ActiveWindow.WindowState = -4143
$excel.ActiveWindow.Top = 0
$excel.ActiveWindow.Left = 0
$excel.ActiveWindow.Width = 28800
$excel.ActiveWindow.Height = 11835
